$wb = $excel.ActiveWorkbook

# --- Sheet "model" (sheet1): insert "status" row after header, append "delivery_name" row ---
$model = $wb.Worksheets.Item("model")

# Insert a new row 2 (pushes existing rows down), copying format from the row above (the header row)
$model.Rows.Item(2).Insert(-4121)  # xlShiftDown; format will come from row above by default behavior in Excel COM? set explicitly below
$model.Range("A2").Value = "string"
$model.Range("B2").Value = "status"

# Append a new "delivery_name" row at the end (row 11), matching style of item_description/is_active/is_simple rows
$model.Range("A11").Value = "string"
$model.Range("B11").Value = "delivery_name"
$model.Range("A10:B10").Copy()
$model.Range("A11:B11").PasteSpecial(-4122) # xlPasteFormats

$model.Range("B6").Select()

# --- Sheet "settings" (sheet2): rename display.title -> display.title.text ---
$settings = $wb.Worksheets.Item("settings")
$settings.Range("C1").Value = "display.title.text"
$settings.Range("C1").Select()

# --- Activate settings sheet as the active tab ---
$settings.Activate()

# --- Update workbook view window position ---
$wb.Windows.Item(1).Left = 10760
$wb.Windows.Item(1).Top = 1900
